$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-01 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-02 Sunday", 2) | Out-Null
$d.Content.Find.Execute("550×9=", $true, $false, $false, $false, $false, $true, 1, $false, "811×4=", 2) | Out-Null
$d.Content.Find.Execute("261×6=", $true, $false, $false, $false, $false, $true, 1, $false, "870×9=", 2) | Out-Null
$d.Content.Find.Execute("939×7=", $true, $false, $false, $false, $false, $true, 1, $false, "153×4=", 2) | Out-Null
$d.Content.Find.Execute("679×7=", $true, $false, $false, $false, $false, $true, 1, $false, "585×7=", 2) | Out-Null
$d.Content.Find.Execute("923×5=", $true, $false, $false, $false, $false, $true, 1, $false, "481×9=", 2) | Out-Null
$d.Content.Find.Execute("719×3=", $true, $false, $false, $false, $false, $true, 1, $false, "711×5=", 2) | Out-Null
$d.Content.Find.Execute("987×5=", $true, $false, $false, $false, $false, $true, 1, $false, "169×6=", 2) | Out-Null
$d.Content.Find.Execute("263×2=", $true, $false, $false, $false, $false, $true, 1, $false, "330×3=", 2) | Out-Null
$d.Content.Find.Execute("362×2=", $true, $false, $false, $false, $false, $true, 1, $false, "305×4=", 2) | Out-Null
$d.Content.Find.Execute("553×3=", $true, $false, $false, $false, $false, $true, 1, $false, "281×2=", 2) | Out-Null
$d.Content.Find.Execute("317×8=", $true, $false, $false, $false, $false, $true, 1, $false, "204×5=", 2) | Out-Null
$d.Content.Find.Execute("176×9=", $true, $false, $false, $false, $false, $true, 1, $false, "567×3=", 2) | Out-Null
$d.Content.Find.Execute("640×9=", $true, $false, $false, $false, $false, $true, 1, $false, "400×6=", 2) | Out-Null
$d.Content.Find.Execute("729×7=", $true, $false, $false, $false, $false, $true, 1, $false, "666×2=", 2) | Out-Null
$d.Content.Find.Execute("966×5=", $true, $false, $false, $false, $false, $true, 1, $false, "382×2=", 2) | Out-Null
$d.Content.Find.Execute("726×3=", $true, $false, $false, $false, $false, $true, 1, $false, "605×7=", 2) | Out-Null
$d.Content.Find.Execute("154×7=", $true, $false, $false, $false, $false, $true, 1, $false, "105×2=", 2) | Out-Null
$d.Content.Find.Execute("830×2=", $true, $false, $false, $false, $false, $true, 1, $false, "598×6=", 2) | Out-Null
$d.Content.Find.Execute("851×7=", $true, $false, $false, $false, $false, $true, 1, $false, "300×3=", 2) | Out-Null
$d.Content.Find.Execute("358×3=", $true, $false, $false, $false, $false, $true, 1, $false, "834×5=", 2) | Out-Null
$d.Content.Find.Execute("849×2=", $true, $false, $false, $false, $false, $true, 1, $false, "745×2=", 2) | Out-Null
$d.Content.Find.Execute("238×3=", $true, $false, $false, $false, $false, $true, 1, $false, "611×3=", 2) | Out-Null
$d.Content.Find.Execute("935×2=", $true, $false, $false, $false, $false, $true, 1, $false, "637×4=", 2) | Out-Null
$d.Content.Find.Execute("271×5=", $true, $false, $false, $false, $false, $true, 1, $false, "564×4=", 2) | Out-Null
$d.Content.Find.Execute("328×7=", $true, $false, $false, $false, $false, $true, 1, $false, "965×5=", 2) | Out-Null
